# Insert a new weekly record above row 129. This pushes the existing
# rows 129-154 down to 130-155 (so the row that falls off the bottom,
# the former row 154, becomes the new row 155), and leaves a blank
# row 129 ready for the new data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(129).Insert()

# The newly inserted row 129 is blank; populate it by cloning the
# (now shifted-down) row 130 -- which still holds the values that used
# to live in row 129 -- for every column except the date (D) and the
# volume (J), which get the new week's figures.
for ($col = 1; $col -le 18; $col++) {
    if ($col -ne 4 -and $col -ne 10) {
        $ws.Cells.Item(129, $col).Value = $ws.Cells.Item(130, $col).Value2
    }
}

$ws.Cells.Item(129, 4).Value = 44504
$ws.Cells.Item(129, 10).Value = 25
